# Apply the "new version with timestamp" edit:
#  1) Insert a new product row ("FAWAR FRUIT 6 SACHETS") right before the
#     "GINKGO BILOBA 30 CAPS." row, renumbering the rank column (A) for every
#     row that shifts down, and recomputing the running total.
#  2) Bump the generated timestamp string shown in the report footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert the new row just above "GINKGO BILOBA 30 CAPS." (row 18).
#    Inserting here pushes every row below it (old rows 18-36) down by one,
#    carrying values, merges (mostly) and row heights along with them.
# ---------------------------------------------------------------------
$newRow = 18
$ws.Rows($newRow.ToString() + ":" + $newRow.ToString()).Insert()

# Re-create the merged cell groups for the freshly inserted row (Insert()
# shifts existing merges down but does not clone them for the blank row).
$ws.Range("A" + $newRow + ":B" + $newRow).Merge()
$ws.Range("C" + $newRow + ":G" + $newRow).Merge()
$ws.Range("H" + $newRow + ":K" + $newRow).Merge()
$ws.Range("L" + $newRow + ":M" + $newRow).Merge()
$ws.Range("N" + $newRow + ":O" + $newRow).Merge()

# Match the row height used by the row before the edit (this stripe keeps
# using 24.75 at this position, same as the rest of the table pattern).
$ws.Rows($newRow).RowHeight = 24.75

# Fill in the new item's data (rank gets fixed up in the renumber pass below).
$ws.Range("A" + $newRow).Value = 12
$ws.Range("C" + $newRow).Value = "FAWAR FRUIT 6 SACHETS"
$ws.Range("H" + $newRow).Value = "8:3"
$ws.Range("L" + $newRow).Value = "1"
$ws.Range("N" + $newRow).Value = "24.00"
$ws.Range("P" + $newRow).Value = "3.8400"
$ws.Range("Q" + $newRow).Value = "0:1"

# ---------------------------------------------------------------------
# 2. Renumber column A (the rank / item number) for every row from the new
#    row through the last product row so the sequence stays 1..29.
# ---------------------------------------------------------------------
$lastProductRow = 35
for ($r = $newRow + 1; $r -le $lastProductRow; $r++) {
    $ws.Range("A" + $r).Value = $r - 6
}

# ---------------------------------------------------------------------
# 3. Recompute the running total shown right below the table.
# ---------------------------------------------------------------------
$totalRow = $lastProductRow + 1
$total = 0
for ($r = 7; $r -le $lastProductRow; $r++) {
    $total = $total + $ws.Range("P" + $r).Value2
}
$ws.Range("P" + $totalRow).Value = $total

# ---------------------------------------------------------------------
# 4. Bump the "generated at" timestamp shown in the footer.
# ---------------------------------------------------------------------
$footerRow = $totalRow + 1
$stampCell = $ws.Range("A" + $footerRow)
$stampCell.Value = "Wednesday, 3 September, 2025 12:34 PM"
